$wb = $excel.ActiveWorkbook

$ws0D   = $wb.Worksheets.Item("0D")
$wsYear = $wb.Worksheets.Item("year_Vecteurs")
$wsProd = $wb.Worksheets.Item("Production_system")

# --- Sheet "0D": drop the two "=1/1000" formulas, replace with the literal value 1,
#     and move the selection to D5 (no longer the tab-selected sheet). ---
$ws0D.Range("B6").Value = 1
$ws0D.Range("B7").Value = 1
$ws0D.Range("D5").Select()

# --- Sheet "Production_system": selection only -> B8 ---
$wsProd.Range("B8").Select()

# --- Sheet "year_Vecteurs": becomes the active/tab-selected sheet, selection -> G12,
#     bold A2:B13, and rescale the "C" values (divide by 1000). ---
$wsYear.Range("A2:B13").Style = "Normal"
$wsYear.Range("A2:B13").Font.Bold = $true

$wsYear.Range("C2").Value = 0.06
$wsYear.Range("C3").Value = 0.006
$wsYear.Range("C6").Value = 0.2
$wsYear.Range("C8").Value = 0.03
$wsYear.Range("C9").Value = 0.006
$wsYear.Range("C12").Value = 0.06

# Activate this sheet last so it ends up as the workbook's active tab, and select
# its range last so the Select() call doesn't re-activate another sheet afterwards.
$wsYear.Activate()
$wsYear.Range("G12").Select()
